$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6: a new tag recapture record
$ws.Range("A6").Value = "NA"
$ws.Range("B6").Value = "NA"
$ws.Range("C6").Value = 43707
$ws.Range("C6").NumberFormat = "d-mmm"
$ws.Range("M6").Value = "no tag number; acoustic tag missing"
$ws.Range("D6").Value = "Lower Harrison Reiver"
$ws.Range("E6").Value = 49.251817000000003
$ws.Range("F6").Value = -121.93721499999999
$ws.Range("G6").Value = "gillnet"
$ws.Range("H6").Value = "Sts'ailes"
$ws.Range("I6").Value = "Kim Charlie"
$ws.Range("J6").Value = 16047962116
$ws.Range("K6").Value = "NA"
$ws.Range("L6").Value = "4690 Salish Way, Agassiz BC, V0M1A1"

# Row 2: fill in previously "NA" tag-recapture details now available
$ws.Range("C2").Value = 43705
$ws.Range("C2").NumberFormat = "d-mmm"
$ws.Range("D2").Value = "Hood Canal"
$ws.Range("E2").Value = 47.407499
$ws.Range("F2").Value = -123.136264

# Row 5: correct the group/contact info for this Squaxin Island fish
$ws.Range("H5").Value = "Squaxin"
$ws.Range("I5").Value = "Mike Foster"

Write-Output "done"
